$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the outlier-flag columns
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match the formatting already used by the other header cells (bold, centered, bordered)
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Boolean outlier flags for each data row (rows 2-21)
$values = @(
    @($false, $false, $true),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $true),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($true,  $true,  $true),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($true,  $true,  $true),
    @($false, $false, $false),
    @($false, $false, $false)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i][0]
    $ws.Cells.Item($row, 7).Value = $values[$i][1]
    $ws.Cells.Item($row, 8).Value = $values[$i][2]
}
